{"js": "const pairs = [\n  [\"2025-02-20 Thursday\", \"2025-02-21 Friday\"],\n  [\"75\u00f74=18, 3\", \"26\u00f78=3, 2\"],\n  [\"39\u00f77=5, 4\", \"52\u00f73=17, 1\"],\n  [\"49\u00f73=16, 1\", \"12\u00f76=2, 0\"],\n  [\"10\u00f77=1, 3\", \"20\u00f79=2, 2\"],\n  [\"17\u00f78=2, 1\", \"51\u00f78=6, 3\"],\n  [\"97\u00f74=24, 1\", \"12\u00f79=1, 3\"],\n  [\"60\u00f72=30, 0\", \"70\u00f75=14, 0\"],\n  [\"24\u00f74=6, 0\", \"24\u00f75=4, 4\"],\n  [\"84\u00f75=16, 4\", \"67\u00f79=7, 4\"],\n  [\"64\u00f73=21, 1\", \"98\u00f79=10, 8\"],\n  [\"98\u00f75=19, 3\", \"34\u00f75=6, 4\"],\n  [\"48\u00f73=16, 0\", \"56\u00f77=8, 0\"],\n  [\"27\u00f75=5, 2\", \"58\u00f78=7, 2\"],\n  [\"50\u00f77=7, 1\", \"65\u00f73=21, 2\"],\n  [\"69\u00f73=23, 0\", \"23\u00f77=3, 2\"],\n  [\"38\u00f76=6, 2\", \"39\u00f74=9, 3\"],\n  [\"26\u00f77=3, 5\", \"33\u00f75=6, 3\"],\n  [\"24\u00f72=12, 0\", \"93\u00f73=31, 0\"],\n  [\"65\u00f76=10, 5\", \"19\u00f75=3, 4\"],\n  [\"30\u00f77=4, 2\", \"59\u00f72=29, 1\"],\n  [\"51\u00f75=10, 1\", \"29\u00f76=4, 5\"],\n  [\"44\u00f75=8, 4\", \"32\u00f74=8, 0\"],\n  [\"44\u00f78=5, 4\", \"42\u00f74=10, 2\"],\n  [\"71\u00f76=11, 5\", \"84\u00f78=10, 4\"],\n  [\"36\u00f77=5, 1\", \"93\u00f74=23, 1\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-02-20 Thursday\", \"2025-02-21 Friday\")\n    ,@(\"75\u00f74=18, 3\", \"26\u00f78=3, 2\")\n    ,@(\"39\u00f77=5, 4\", \"52\u00f73=17, 1\")\n    ,@(\"49\u00f73=16, 1\", \"12\u00f76=2, 0\")\n    ,@(\"10\u00f77=1, 3\", \"20\u00f79=2, 2\")\n    ,@(\"17\u00f78=2, 1\", \"51\u00f78=6, 3\")\n    ,@(\"97\u00f74=24, 1\", \"12\u00f79=1, 3\")\n    ,@(\"60\u00f72=30, 0\", \"70\u00f75=14, 0\")\n    ,@(\"24\u00f74=6, 0\", \"24\u00f75=4, 4\")\n    ,@(\"84\u00f75=16, 4\", \"67\u00f79=7, 4\")\n    ,@(\"64\u00f73=21, 1\", \"98\u00f79=10, 8\")\n    ,@(\"98\u00f75=19, 3\", \"34\u00f75=6, 4\")\n    ,@(\"48\u00f73=16, 0\", \"56\u00f77=8, 0\")\n    ,@(\"27\u00f75=5, 2\", \"58\u00f78=7, 2\")\n    ,@(\"50\u00f77=7, 1\", \"65\u00f73=21, 2\")\n    ,@(\"69\u00f73=23, 0\", \"23\u00f77=3, 2\")\n    ,@(\"38\u00f76=6, 2\", \"39\u00f74=9, 3\")\n    ,@(\"26\u00f77=3, 5\", \"33\u00f75=6, 3\")\n    ,@(\"24\u00f72=12, 0\", \"93\u00f73=31, 0\")\n    ,@(\"65\u00f76=10, 5\", \"19\u00f75=3, 4\")\n    ,@(\"30\u00f77=4, 2\", \"59\u00f72=29, 1\")\n    ,@(\"51\u00f75=10, 1\", \"29\u00f76=4, 5\")\n    ,@(\"44\u00f75=8, 4\", \"32\u00f74=8, 0\")\n    ,@(\"44\u00f78=5, 4\", \"42\u00f74=10, 2\")\n    ,@(\"71\u00f76=11, 5\", \"84\u00f78=10, 4\")\n    ,@(\"36\u00f77=5, 1\", \"93\u00f74=23, 1\")\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p[0]\n    $find.Replacement.Text = $p[1]\n    $find.Execute($p[0], $false, $false, $false, $false, $false, $true, 1, $false, $p[1], 2)\n}\n"}
